# Weekly update: insert a new price report row for "Feria Lagunitas de
# Puerto Montt - Coliflor" ahead of the existing row 242, shifting every
# subsequent row down by one (old row 242 -> 243, ..., old row 262 -> 263).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 242:262 down to 243:263, leaving a blank row 242 to fill in.
$ws.Rows.Item(242).Insert()

$ws.Cells.Item(242, 1).Value  = 4
$ws.Cells.Item(242, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(242, 3).Value  = "Los Lagos"
$ws.Cells.Item(242, 4).Value  = 44585
$ws.Cells.Item(242, 5).Value  = 10
$ws.Cells.Item(242, 6).Value  = 100112008
$ws.Cells.Item(242, 7).Value  = "Coliflor"
$ws.Cells.Item(242, 8).Value  = "Sin especificar"
$ws.Cells.Item(242, 9).Value  = "Primera"
$ws.Cells.Item(242, 10).Value = 250
$ws.Cells.Item(242, 11).Value = 1500
$ws.Cells.Item(242, 12).Value = 1500
$ws.Cells.Item(242, 13).Value = 1500
$ws.Cells.Item(242, 14).Value = "`$/unidad"
$ws.Cells.Item(242, 15).Value = "Región Metropolitana"
$ws.Cells.Item(242, 16).Value = 1500
$ws.Cells.Item(242, 17).Value = 1
$ws.Cells.Item(242, 18).Value = "Hortaliza"
